# docs/ValueSet-AllAppIdCodes.xlsx — update MIN / MAX values to align with
# MHV-17222 (bump FHIR IG "Version" and regenerate the publication "Date"
# on the Metadata sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 3: Property "Version" -> Value column (B3)
$ws.Range("B3").Value = "0.2.10-beta"

# Row 8: Property "Date" -> Value column (B8)
$ws.Range("B8").Value = "2023-12-06T12:46:33-06:00"
